$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new data row at row 85 (pushing the old "totals" row 85 -> 86
#     and the old "footer" row 86 -> 87), mirroring the row that was added
#     for the new item "مناديل FINE" in the source pharmacy report.

$ws.Rows("85:85").Insert()

# Copy the formatting (styles/merges) of the row above (the previous last
# item row, row 84) into the freshly inserted row 85, then set its own
# height explicitly (Excel recalculates wrap-height per row; the new row
# ends up at 25.5pt here).
$ws.Range("A84:N84").Copy($ws.Range("A85:N85"))
$ws.Rows("85:85").RowHeight = 25.5

# New item values for row 85
$ws.Range("A85").Value2 = 82
$ws.Range("B85").Value2 = "مناديل FINE"
$ws.Range("H85").Value2 = "15:0"
$ws.Range("L85").Value2 = 30
$ws.Range("N85").Value2 = "1:0"

# The old "totals" row (now shifted down to row 86) picks up the new
# item's contribution to the K total (5119.27 -> 5149.27).
$ws.Range("K86").Value2 = 5149.2700000000004
